$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19, shifting existing rows 19-31 down to 20-32.
$ws.Rows.Item(19).Insert()

# B19 should have the plain bordered checklist style, same as B20.
$ws.Range("B20").Copy()
$ws.Range("B19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# A19 should carry the "done" highlight style, matching the other completed
# (green-filled) checkbox cells such as A18/A21.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's content.
$ws.Range("B19").Value = "add user feedback for clicks, successful set, total sets found"

# Select the new active cell as recorded post-edit.
$ws.Range("B23").Select()
